$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LinksLogin")

# Widen column A from 18.21875 to 29.21875 characters.
$ws.Columns.Item(1).ColumnWidth = 29.21875

# Update the stored password value in B2.
$ws.Cells.Item(2, 2).Value = 123456789
